# "Interaccion con el calendario" - the user picked a date from the date
# picker on the Recruitment sheet's application-date field, which (a) fixed
# the header label typo, (b) rewrote the date value in DD-MM-YYYY form as
# produced by the calendar widget, and (c) left the active selection on the
# calendar's target cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recruitment")

# Fix the column header (DATA_OF_APPLICATION -> DATE_OF_APPLICATION)
$ws.Range("I1").Value = "DATE_OF_APPLICATION"

# New date chosen via the calendar control (2024-22-07 -> 17-07-2024)
$ws.Range("I2").Value = "17-07-2024"

# The calendar interaction leaves the selection on the next row's cell
[void]$ws.Range("I6").Select()
